# Apply the "Hypotheses_paper_pulp_1D" update:
#  - translate the 9 pulping-method scenario names to French
#  - add a new "seasonal_efficiency" column to year_Vecteurs
#  - tidy up formatting (add thin box borders around the used ranges
#    that were missing them, and fix the partial border on the
#    retrofit_Transition header)
#  - restore per-sheet selections / active sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Rename the 9 pulping-method labels everywhere they are used.
#    (Production_system, Production_system_year and
#    retrofit_Transition all reference these shared strings.)
# ---------------------------------------------------------------
$renames = @{
    "Mechanical pulping + Fossil fired PM" = "Pulping mecanique + MP fossile";
    "Chemical pulping + Fossil fired PM"   = "Pulping chimique + MP fossile";
    "DES pulping + Fossil fired PM"        = "Pulping via DES + MP fossile";
    "Mechanical pulping + Electric PM"     = "Pulping mecanique + MP electrique";
    "Chemical pulping + Electric PM"       = "Pulping chimique + MP electrique";
    "DES pulping + Electric PM"            = "Pulping via DES + MP electrique";
    "Mechanical pulping + Biomass PM"      = "Pulping mecanique + MP bois";
    "Chemical pulping + Biomass PM"        = "Pulping chimique + MP bois";
    "DES pulping + Biomass PM"             = "Pulping via DES + MP bois";
}

foreach ($ws in $wb.Worksheets) {
    foreach ($oldName in $renames.Keys) {
        $ws.Cells.Replace($oldName, $renames[$oldName]) | Out-Null
    }
}

# ---------------------------------------------------------------
# 2) year_Vecteurs: add the new "seasonal_efficiency" column (E)
# ---------------------------------------------------------------
$wsYV = $wb.Worksheets.Item("year_Vecteurs")

$wsYV.Range("E1").Value = "seasonal_efficiency"

$wsYV.Range("E2").Value  = 0.64724947431429825
$wsYV.Range("E3").Value  = 0.81021150558322341
$wsYV.Range("E4").Value  = 0.69400671890583965
$wsYV.Range("E5").Value  = 0.61170901977899528
$wsYV.Range("E6").Value  = 0.63115005598290974
$wsYV.Range("E7").Value  = 0.81021150558322341
$wsYV.Range("E8").Value  = 1.0351908199407167
$wsYV.Range("E9").Value  = 0.91962472240630677
$wsYV.Range("E10").Value = 0.85376273680792314
$wsYV.Range("E11").Value = 0.69348451972740355
$wsYV.Range("E12").Value = 0.73620708377192978
$wsYV.Range("E13").Value = 0.91962472240630677

# Tidy formatting: add a thin box border around the whole used range
# (also normalises the highlighted hydrogen rows to keep their fill).
$wsYV.Range("A1:E13").Borders.LineStyle = 1

# ---------------------------------------------------------------
# 3) Production_system: add a thin box border around the used range
# ---------------------------------------------------------------
$wsPS = $wb.Worksheets.Item("Production_system")
$wsPS.Range("A1:C10").Borders.LineStyle = 1

# ---------------------------------------------------------------
# 4) retrofit_Transition: add a thin box border around the used
#    range (this also repairs B1's partial border).
# ---------------------------------------------------------------
$wsRT = $wb.Worksheets.Item("retrofit_Transition")
$wsRT.Range("A1:K10").Borders.LineStyle = 1

# ---------------------------------------------------------------
# 5) Restore per-sheet selections.
# ---------------------------------------------------------------
$wsPS.Activate() | Out-Null
$wsPS.Range("A1:C10").Select() | Out-Null

$wsPY = $wb.Worksheets.Item("Production_system_year")
$wsPY.Activate() | Out-Null
$wsPY.Range("A24").Select() | Out-Null

$wsYV.Activate() | Out-Null
$wsYV.Range("E23").Select() | Out-Null

$wsRT.Activate() | Out-Null
$wsRT.Range("A18").Select() | Out-Null

# "0D" stays the active / tab-selected sheet, so activate it last.
$ws0D = $wb.Worksheets.Item("0D")
$ws0D.Activate() | Out-Null
$ws0D.Range("B32").Select() | Out-Null
